# Shift the weekly Perejil price rows down by one: a brand-new observation is
# inserted at row 83, pushing every existing row (83..165) down by one so the
# previous row 165 becomes row 166. Only columns D (Fecha), J (Volumen),
# K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado) and
# P (Precio $/Kg) vary row to row; the remaining columns are constant for
# every record in this sheet/subset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 83
$lastRow  = 165
$newLastRow = $lastRow + 1

# --- 1. Capture the "old" values for the columns that move, before touching
#        anything (indexed by source row number). ---
$oldD = @{}
$oldJ = @{}
$oldK = @{}
$oldL = @{}
$oldM = @{}
$oldP = @{}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $oldD[$r] = $ws.Range("D$r").Value2
    $oldJ[$r] = $ws.Range("J$r").Value2
    $oldK[$r] = $ws.Range("K$r").Value2
    $oldL[$r] = $ws.Range("L$r").Value2
    $oldM[$r] = $ws.Range("M$r").Value2
    $oldP[$r] = $ws.Range("P$r").Value2
}

# --- 2. Row 166 is brand new: populate the columns that are constant across
#        the whole block by copying them from row 165 (any source row would
#        do, they are all identical). NOTE: reads must use .Value2 — .Value
#        does not resolve to the underlying scalar in this COM shim. ---
$ws.Range("A$newLastRow").Value = $ws.Range("A$lastRow").Value2
$ws.Range("B$newLastRow").Value = $ws.Range("B$lastRow").Value2
$ws.Range("C$newLastRow").Value = $ws.Range("C$lastRow").Value2
$ws.Range("E$newLastRow").Value = $ws.Range("E$lastRow").Value2
$ws.Range("F$newLastRow").Value = $ws.Range("F$lastRow").Value2
$ws.Range("G$newLastRow").Value = $ws.Range("G$lastRow").Value2
$ws.Range("H$newLastRow").Value = $ws.Range("H$lastRow").Value2
$ws.Range("I$newLastRow").Value = $ws.Range("I$lastRow").Value2
$ws.Range("N$newLastRow").Value = $ws.Range("N$lastRow").Value2
$ws.Range("O$newLastRow").Value = $ws.Range("O$lastRow").Value2
$ws.Range("Q$newLastRow").Value = $ws.Range("Q$lastRow").Value2
$ws.Range("R$newLastRow").Value = $ws.Range("R$lastRow").Value2

# Give the new D cell the same date-style number format as the others.
$ws.Range("D$newLastRow").NumberFormat = $ws.Range("D$lastRow").NumberFormat

# --- 3. Shift D/J/K/L/M/P down by one row: row r (84..166) gets what used to
#        be in row r-1. Walk from the bottom up so we never overwrite a
#        source value before it has been used (values were captured in step
#        1 anyway, so order is not strictly required, but keep it safe). ---
for ($r = $newLastRow; $r -ge ($firstRow + 1); $r--) {
    $src = $r - 1
    $ws.Range("D$r").Value = $oldD[$src]
    $ws.Range("J$r").Value = $oldJ[$src]
    $ws.Range("K$r").Value = $oldK[$src]
    $ws.Range("L$r").Value = $oldL[$src]
    $ws.Range("M$r").Value = $oldM[$src]
    $ws.Range("P$r").Value = $oldP[$src]
}

# --- 4. Row 83 becomes the new observation: new date + new volume, while
#        the min/max/weighted-avg/$-per-kg prices are unchanged. ---
$ws.Range("D$firstRow").Value = 44810
$ws.Range("J$firstRow").Value = 2400
$ws.Range("K$firstRow").Value = $oldK[$firstRow]
$ws.Range("L$firstRow").Value = $oldL[$firstRow]
$ws.Range("M$firstRow").Value = $oldM[$firstRow]
$ws.Range("P$firstRow").Value = $oldP[$firstRow]
